# "fix hearing value sets"
#
# Target layout (tab order):
#   1. Metadata                               (unchanged, except the Date value)
#   2. Include from Hearing Loss - T          (renamed from "Include ValueSets", 2 rows appended)
#   3. Include from Hearing Loss - D          (renamed from "Include ValueSets 2", 2 rows appended)
#   4. Include ValueSets                      (was "Include ValueSets 3")
#   5. Include ValueSets 2                    (was "Include ValueSets 4")
#   6. Include ValueSets 3                    (was "Include ValueSets 5")
#   7. Include ValueSets 4                    (new)
#   8. Include ValueSets 5                    (new)
#
# sheetId/rId in the target file are perfectly sequential (1..8) in tab order,
# which only happens if the "Include ValueSets*" sheets are rebuilt from
# scratch after the Metadata sheet (Excel hands out sheetId = current sheet
# count + 1 at creation time) - so we delete the five original include
# sheets and recreate all seven in the right order/position.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Metadata sheet: bump the "Date" property value (row 8, column B)
# ---------------------------------------------------------------------
$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B8").Value = "2022-05-04T12:01:11-04:00"

# ---------------------------------------------------------------------
# 2. Remember header/body style source cells from the Metadata sheet
#    style 1 = bold header row, style 2 = plain body row (with borders)
# ---------------------------------------------------------------------
$headerStyleSrc = $metadata.Range("A1")
$bodyStyleSrc   = $metadata.Range("A2")

# ---------------------------------------------------------------------
# 3. Drop the five original "Include ValueSets*" sheets
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Include ValueSets 5").Delete() | Out-Null
$wb.Worksheets.Item("Include ValueSets 4").Delete() | Out-Null
$wb.Worksheets.Item("Include ValueSets 3").Delete() | Out-Null
$wb.Worksheets.Item("Include ValueSets 2").Delete() | Out-Null
$wb.Worksheets.Item("Include ValueSets").Delete() | Out-Null

# ---------------------------------------------------------------------
# Helper-ish inline blocks: build each replacement sheet right after the
# previous one so the tab order comes out exactly right.
# ---------------------------------------------------------------------

# ---- Sheet: "Include from Hearing Loss - T" -----------------------------
$prev = $metadata
$s2 = $wb.Worksheets.Add($null, $prev)
$s2.Name = "Include from Hearing Loss - T"
$s2.Columns.Item(1).ColumnWidth = 30.703125
$s2.Columns.Item(2).ColumnWidth = 50.703125

$s2.Range("A1").Value = "Codes"
$s2.Range("A2").Value = "All codes"
$s2.Range("A3").Value = ""
$s2.Range("B3").Value = ""
$s2.Range("A4").Value = "System URI"
$s2.Range("B4").Value = "http://hl7.org/fhir/us/pacio-splasch/CodeSystem/HearingLossTypeCS"

$headerStyleSrc.Copy() | Out-Null
$s2.Range("A1").PasteSpecial($xlPasteFormats) | Out-Null
$bodyStyleSrc.Copy() | Out-Null
$s2.Range("A2:B4").PasteSpecial($xlPasteFormats) | Out-Null

# ---- Sheet: "Include from Hearing Loss - D" -----------------------------
$prev = $s2
$s3 = $wb.Worksheets.Add($null, $prev)
$s3.Name = "Include from Hearing Loss - D"
$s3.Columns.Item(1).ColumnWidth = 30.703125
$s3.Columns.Item(2).ColumnWidth = 50.703125

$s3.Range("A1").Value = "Codes"
$s3.Range("A2").Value = "All codes"
$s3.Range("A3").Value = ""
$s3.Range("B3").Value = ""
$s3.Range("A4").Value = "System URI"
$s3.Range("B4").Value = "http://hl7.org/fhir/us/pacio-splasch/CodeSystem/HearingLossDegreeCS"

$headerStyleSrc.Copy() | Out-Null
$s3.Range("A1").PasteSpecial($xlPasteFormats) | Out-Null
$bodyStyleSrc.Copy() | Out-Null
$s3.Range("A2:B4").PasteSpecial($xlPasteFormats) | Out-Null

# ---- Sheet: "Include ValueSets" -----------------------------------------
$prev = $s3
$s4 = $wb.Worksheets.Add($null, $prev)
$s4.Name = "Include ValueSets"
$s4.Columns.Item(1).ColumnWidth = 30.703125
$s4.Columns.Item(2).ColumnWidth = 50.703125

$s4.Range("A1").Value = "ValueSet URL"
$s4.Range("A2").Value = "http://hl7.org/fhir/us/pacio-splasch/ValueSet/HearingAbilityToHearDuringAssessmentPeriodVS"

$headerStyleSrc.Copy() | Out-Null
$s4.Range("A1").PasteSpecial($xlPasteFormats) | Out-Null
$bodyStyleSrc.Copy() | Out-Null
$s4.Range("A2").PasteSpecial($xlPasteFormats) | Out-Null

# ---- Sheet: "Include ValueSets 2" ----------------------------------------
$prev = $s4
$s5 = $wb.Worksheets.Add($null, $prev)
$s5.Name = "Include ValueSets 2"
$s5.Columns.Item(1).ColumnWidth = 30.703125
$s5.Columns.Item(2).ColumnWidth = 50.703125

$s5.Range("A1").Value = "ValueSet URL"
$s5.Range("A2").Value = "http://hl7.org/fhir/us/pacio-splasch/ValueSet/YesOrNoVS"

$headerStyleSrc.Copy() | Out-Null
$s5.Range("A1").PasteSpecial($xlPasteFormats) | Out-Null
$bodyStyleSrc.Copy() | Out-Null
$s5.Range("A2").PasteSpecial($xlPasteFormats) | Out-Null

# ---- Sheet: "Include ValueSets 3" ----------------------------------------
$prev = $s5
$s6 = $wb.Worksheets.Add($null, $prev)
$s6.Name = "Include ValueSets 3"
$s6.Columns.Item(1).ColumnWidth = 30.703125
$s6.Columns.Item(2).ColumnWidth = 50.703125

$s6.Range("A1").Value = "ValueSet URL"
$s6.Range("A2").Value = "http://hl7.org/fhir/us/pacio-splasch/ValueSet/CmsFasiAssistiveDevicesVS"

$headerStyleSrc.Copy() | Out-Null
$s6.Range("A1").PasteSpecial($xlPasteFormats) | Out-Null
$bodyStyleSrc.Copy() | Out-Null
$s6.Range("A2").PasteSpecial($xlPasteFormats) | Out-Null

# ---- Sheet: "Include ValueSets 4" (new) ----------------------------------
$prev = $s6
$s7 = $wb.Worksheets.Add($null, $prev)
$s7.Name = "Include ValueSets 4"
$s7.Columns.Item(1).ColumnWidth = 30.703125
$s7.Columns.Item(2).ColumnWidth = 50.703125

$s7.Range("A1").Value = "ValueSet URL"
$s7.Range("A2").Value = "http://hl7.org/fhir/us/pacio-splasch/ValueSet/LeftRightBothVS"

$headerStyleSrc.Copy() | Out-Null
$s7.Range("A1").PasteSpecial($xlPasteFormats) | Out-Null
$bodyStyleSrc.Copy() | Out-Null
$s7.Range("A2").PasteSpecial($xlPasteFormats) | Out-Null

# ---- Sheet: "Include ValueSets 5" (new) ----------------------------------
$prev = $s7
$s8 = $wb.Worksheets.Add($null, $prev)
$s8.Name = "Include ValueSets 5"
$s8.Columns.Item(1).ColumnWidth = 30.703125
$s8.Columns.Item(2).ColumnWidth = 50.703125

$s8.Range("A1").Value = "ValueSet URL"
$s8.Range("A2").Value = "http://hl7.org/fhir/us/pacio-splasch/ValueSet/HearBetterInOneEarVS"

$headerStyleSrc.Copy() | Out-Null
$s8.Range("A1").PasteSpecial($xlPasteFormats) | Out-Null
$bodyStyleSrc.Copy() | Out-Null
$s8.Range("A2").PasteSpecial($xlPasteFormats) | Out-Null

# ---------------------------------------------------------------------
# Restore original active sheet/selection (Metadata, tab 0)
# ---------------------------------------------------------------------
$metadata.Activate()
$metadata.Range("A1").Select() | Out-Null
